# Applies the data refresh captured by the diff:
#  - Normalizes the "开始时间" (start date) column (B) from dot-separated
#    (2024.03.02) to dash-separated (2024-03-02) format.
#  - Bumps a handful of "想去人数" (interest count) values in column F.
# The same data lives (duplicated) on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Column B (start date) updates: row number -> new date string
$dateUpdates = @{
    2  = "2024-03-02"
    3  = "2024-03-16"
    4  = "2024-03-17"
    5  = "2024-03-23"
    6  = "2024-04-04"
    7  = "2024-04-04"
    8  = "2024-04-04"
    9  = "2024-04-04"
    10 = "2024-04-21"
    11 = "2024-05-03"
    12 = "2024-05-18"
}

# Column F (想去人数 / interest count) updates: row number -> new numeric value
$countUpdates = @{
    3 = 2635
    5 = 87
    6 = 6599
    7 = 447
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $dateUpdates.Keys) {
        $cell = $ws.Cells.Item($row, 2)
        # Force the cell to stay plain text (matching the inline string in
        # the source file) instead of letting Excel auto-convert the
        # dash-separated value into a date serial number.
        $cell.NumberFormat = "@"
        $cell.Value = $dateUpdates[$row]
    }

    foreach ($row in $countUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $countUpdates[$row]
    }
}
